$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds text-formatted numeric-looking strings (prices).
# Prefix with a literal quote to force text entry (matches how Excel
# itself keeps such values as text), then reset the resulting
# quote-prefix style back to Normal so no visible style/format changes.

# Row 2
$ws.Range('D2').Value = "'54.693.83"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.41%  '

# Row 3
$ws.Range('D3').Value = "'2.322.53"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -6.91%  '

# Row 4
$ws.Range('E4').Value = '  +0.08%  '

# Row 5
$ws.Range('D5').Value = "'470.09"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.47%  '

# Row 6
$ws.Range('D6').Value = "'142.95"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.16%  '

# Row 7
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').Value = "'1.00"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.40%  '

# Row 8
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').Value = "'0.611"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +18.70%  '

# Row 9
$ws.Range('D9').Value = "'2.327.80"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -7.22%  '

# Row 10
$ws.Range('D10').Value = "'0.0948"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.48%  '

# Row 11
$ws.Range('D11').Value = "'5.41"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -7.24%  '

# Row 12
$ws.Range('D12').Value = "'0.323"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.45%  '

# Row 13
$ws.Range('E13').Value = '  +0.60%  '

# Row 14
$ws.Range('D14').Value = "'2.740.82"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.56%  '

# Row 15
$ws.Range('D15').Value = "'54.874.12"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.89%  '

# Row 16
$ws.Range('D16').Value = "'19.77"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -7.40%  '

# Row 17
$ws.Range('D17').Value = "'0.0000128"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -7.24%  '

# Row 18
$ws.Range('D18').Value = "'2.340.41"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -6.77%  '

# Row 19
$ws.Range('D19').Value = "'4.50"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.48%  '

# Row 20
$ws.Range('D20').Value = "'311.67"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.07%  '

# Row 21
$ws.Range('D21').Value = "'9.48"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -7.14%  '

# Row 22
$ws.Range('D22').Value = "'0.997"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.23%  '

# Row 23
$ws.Range('D23').Value = "'5.58"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.55%  '

# Row 24
$ws.Range('D24').Value = "'56.20"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.76%  '

# Row 25
$ws.Range('D25').Value = "'1.00"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.43%  '

# Row 26
$ws.Range('D26').Value = "'0.390"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -6.38%  '

# Row 27
$ws.Range('D27').Value = "'0.150"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -8.15%  '

# Row 28
$ws.Range('D28').Value = "'2.444.35"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.84%  '

# Row 29
$ws.Range('D29').Value = "'6.97"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -10.58%  '

# Row 30
$ws.Range('E30').Value = '  +0.21%  '

# Row 31
$ws.Range('D31').Value = "'0.0₃0730"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -9.15%  '

# Row 32
$ws.Range('D32').Value = "'145.37"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.42%  '

# Row 33
$ws.Range('D33').Value = "'17.90"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.43%  '

# Row 34
$ws.Range('D34').Value = "'1.45"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.09%  '

# Row 35
$ws.Range('D35').Value = "'5.03"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.11%  '

# Row 36
$ws.Range('D36').Value = "'1.08"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.12%  '

# Row 37
$ws.Range('D37').Value = "'3.55"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.26%  '

# Row 38
$ws.Range('D38').Value = "'0.798"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -9.08%  '

# Row 39
$ws.Range('D39').Value = "'0.101"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +9.12%  '

# Row 40
$ws.Range('D40').Value = "'33.57"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.78%  '

# Row 41
$ws.Range('D41').Value = "'1.00"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.59%  '

# Row 42
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = "'3.35"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.56%  '

# Row 43
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = "'1.31"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.77%  '

# Row 44
$ws.Range('D44').Value = "'0.571"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.72%  '

# Row 45
$ws.Range('B45').Value = 'WhiteBITCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D45').Value = "'10.16"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.46%  '

# Row 46
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').Value = "'0.0512"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -8.74%  '

# Row 47
$ws.Range('D47').Value = "'247.47"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.22%  '

# Row 48
$ws.Range('D48').Value = "'0.0218"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.98%  '

# Row 49
$ws.Range('D49').Value = "'4.32"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -10.60%  '

# Row 50
$ws.Range('D50').Value = "'1.774.30"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.51%  '

# Row 51
$ws.Range('D51').Value = "'16.43"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -7.58%  '
